$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-12-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-09 Monday", 2) | Out-Null

# Update the division problems in the table, addressed by (row, col) to avoid
# ambiguity from duplicate cell text (e.g. "22÷6=3, 4" appears twice originally).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "90÷9=10, 0"
$t.Cell(1, 2).Range.Text = "29÷4=7, 1"
$t.Cell(1, 3).Range.Text = "42÷9=4, 6"
$t.Cell(1, 4).Range.Text = "40÷6=6, 4"
$t.Cell(1, 5).Range.Text = "52÷5=10, 2"
$t.Cell(5, 1).Range.Text = "52÷8=6, 4"
$t.Cell(5, 2).Range.Text = "50÷2=25, 0"
$t.Cell(5, 3).Range.Text = "85÷3=28, 1"
$t.Cell(5, 4).Range.Text = "21÷7=3, 0"
$t.Cell(5, 5).Range.Text = "88÷5=17, 3"
$t.Cell(9, 1).Range.Text = "94÷5=18, 4"
$t.Cell(9, 2).Range.Text = "54÷5=10, 4"
$t.Cell(9, 3).Range.Text = "62÷7=8, 6"
$t.Cell(9, 4).Range.Text = "65÷8=8, 1"
$t.Cell(9, 5).Range.Text = "85÷5=17, 0"
$t.Cell(13, 1).Range.Text = "37÷4=9, 1"
$t.Cell(13, 2).Range.Text = "11÷8=1, 3"
$t.Cell(13, 3).Range.Text = "76÷6=12, 4"
$t.Cell(13, 4).Range.Text = "30÷9=3, 3"
$t.Cell(13, 5).Range.Text = "54÷4=13, 2"
$t.Cell(17, 1).Range.Text = "73÷9=8, 1"
$t.Cell(17, 2).Range.Text = "14÷9=1, 5"
$t.Cell(17, 3).Range.Text = "99÷5=19, 4"
$t.Cell(17, 4).Range.Text = "34÷4=8, 2"
$t.Cell(17, 5).Range.Text = "20÷8=2, 4"
